$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("medians")

$ws.Range("C2").Value = 0.279919560267477
$ws.Range("D2").Value = 27.396470319634702
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 5351.1918492887398

$ws.Range("C3").Value = 271.65007443571898
$ws.Range("D3").Value = 28.177374429223701
$ws.Range("E3").Value = 37096.005073127599

$ws.Range("F4").Value = 5348.10569652698

$ws.Range("C5").Value = 0.342140974932367
$ws.Range("D5").Value = 27.686109589041099
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 4405.8679354094602

$ws.Range("C6").Value = 234.22849527491201
$ws.Range("D6").Value = 27.9480867579909
$ws.Range("E6").Value = 22061.897790453

$ws.Range("F7").Value = 4398.7208445469696

$ws.Range("M2").Select()
